# Updates cryptos list values (price/volume, and two rank swaps) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.828.55"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.132.13"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.39"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.77"
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.121.17"
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("E10").Value = "  +6.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.70"
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.95"
$ws.Range("E14").Value = "  +4.14%  "
$ws.Range("E15").Value = "  -1.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.653.54"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.567.22"
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.133.26"
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.09"
$ws.Range("E19").Value = "  -1.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "464.43"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.25"
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.726"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.43"
$ws.Range("E23").Value = "  -1.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.01"
$ws.Range("E24").Value = "  -2.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.27"
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.20"
$ws.Range("E26").Value = "  +1.98%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.07"
$ws.Range("E28").Value = "  +6.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.68"
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.21"
$ws.Range("E31").Value = "  -0.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.00"
$ws.Range("E32").Value = "  +2.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.86"
$ws.Range("E33").Value = "  -0.20%  "
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0869"
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.40"
$ws.Range("E36").Value = "  +2.21%  "
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("E38").Value = "  -3.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.00"
$ws.Range("E39").Value = "  -1.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.32"
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "440.64"
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.70"
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.908.17"
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0370"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.273"
$ws.Range("E45").Value = "  -2.50%  "
$ws.Range("E46").Value = "  -2.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.11"
$ws.Range("E47").Value = "  +2.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.51"
$ws.Range("E48").Value = "  +1.95%  "
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.31"
$ws.Range("E51").Value = "  -0.90%  "
